$wb = $excel.ActiveWorkbook

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4473.184
$ws.Range("I32").Value = 3587.6494
$ws.Range("J32").Value = 11291.8
$ws.Range("K32").Value = 3587.6494
$ws.Range("L32").Value = 11291.8
$ws.Range("M32").Value = -3300.6494
$ws.Range("N32").Value = -11865.8
$ws.Range("H74").Value = 2571.7358
$ws.Range("I74").Value = 690.93024
$ws.Range("K74").Value = 690.93024
$ws.Range("M74").Value = 183.06976
$ws.Range("H77").Value = 2571.7358
$ws.Range("I77").Value = 690.93024
$ws.Range("K77").Value = 3454.6512
$ws.Range("M77").Value = 913.3487999999998
$ws.Range("H119").Value = 40899.5
$ws.Range("J119").Value = 40899.5
$ws.Range("L119").Value = 40899.5
$ws.Range("N119").Value = -50575.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H120").Value = 48000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 48000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 48000
$ws.Range("N120").Value = -57676
$ws.Range("H122").Value = 47780
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 47780
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 47780
$ws.Range("N122").Value = -57580
$ws.Range("H123").Value = 46665.668
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 46665.668
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 46665.668
$ws.Range("N123").Value = -56465.668
$ws.Range("H124").Value = 37195
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 37195
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 37195
$ws.Range("N124").Value = -47015
$ws.Range("H125").Value = 36056.668
$ws.Range("I125").Value = 40000
$ws.Range("J125").Value = 35268
$ws.Range("K125").Value = 40000
$ws.Range("L125").Value = 35268
$ws.Range("M125").Value = -35080
$ws.Range("N125").Value = -45108
$ws.Range("H126").Value = 15445
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 15445
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15445
$ws.Range("N126").Value = -25325
$ws.Range("H127").Value = 40000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 40000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920
$ws.Range("H128").Value = 1050
$ws.Range("I128").Value = 1050
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 3150
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -660
$ws.Range("H129").Value = 45749.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45749.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45749.5
$ws.Range("N129").Value = -55749.5
$ws.Range("H130").Value = 47780
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 47780
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 47780
$ws.Range("N130").Value = -57820
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 40780
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 40780
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 40780
$ws.Range("N132").Value = -50900
$ws.Range("H133").Value = 57000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 57000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 57000
$ws.Range("N133").Value = -67120
$ws.Range("H134").Value = 751.4727
$ws.Range("I134").Value = 637.34
$ws.Range("J134").Value = 1892.8
$ws.Range("K134").Value = 1912.02
$ws.Range("L134").Value = 5678.4
$ws.Range("M134").Value = 622.98
$ws.Range("N134").Value = -10748.4
$ws.Range("H135").Value = 37000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 37000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 37000
$ws.Range("N135").Value = -47140
$ws.Range("H137").Value = 55000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 55000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200
$ws.Range("H138").Value = 50000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 50000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H139").Value = 65000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 65000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280
$ws.Range("H140").Value = 43853.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 43853.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 43853.332
$ws.Range("N140").Value = -54213.332
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1241.6957
$ws.Range("I4").Value = 414.45456
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 1243.36368
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -1131.36368
$ws.Range("N4").Value = -6224

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 37000
$ws.Range("J119").Value = 37000
$ws.Range("L119").Value = 37000
$ws.Range("N119").Value = -46676

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H119").Value = 46500
$ws.Range("J119").Value = 46500
$ws.Range("L119").Value = 46500
$ws.Range("N119").Value = -56176
